# Weekly update: prepend the newest week's Betarraga (Hortaliza) price rows.
# This pushes all existing data rows down by 2 (dimension grows from R229 to
# R231) and fills the newly opened rows 70/71 with the new week's values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new blank rows at row 70, shifting rows 70:229 down to 72:231.
$ws.Rows.Item(70).Resize(2).Insert()

# New row 70 - "Primera" quality for the new week (Fecha serial 44581).
$ws.Range("A70").Value = 8
$ws.Range("B70").Value = "Terminal La Palmera de La Serena"
$ws.Range("C70").Value = "Coquimbo"
$ws.Range("D70").Value = 44581
$ws.Range("E70").Value = 4
$ws.Range("F70").Value = 100114014
$ws.Range("G70").Value = "Betarraga"
$ws.Range("H70").Value = "Sin especificar"
$ws.Range("I70").Value = "Primera"
$ws.Range("J70").Value = 3080
$ws.Range("K70").Value = 450
$ws.Range("L70").Value = 500
$ws.Range("M70").Value = 475
$ws.Range("N70").Value = "`$/paquete 3 unidades"
$ws.Range("O70").Value = "Provincia del Elquí"
$ws.Range("P70").Value = 158
$ws.Range("Q70").Value = 3
$ws.Range("R70").Value = "Hortaliza"

# New row 71 - "Segunda" quality for the new week (Fecha serial 44581).
$ws.Range("A71").Value = 8
$ws.Range("B71").Value = "Terminal La Palmera de La Serena"
$ws.Range("C71").Value = "Coquimbo"
$ws.Range("D71").Value = 44581
$ws.Range("E71").Value = 4
$ws.Range("F71").Value = 100114014
$ws.Range("G71").Value = "Betarraga"
$ws.Range("H71").Value = "Sin especificar"
$ws.Range("I71").Value = "Segunda"
$ws.Range("J71").Value = 1520
$ws.Range("K71").Value = 350
$ws.Range("L71").Value = 400
$ws.Range("M71").Value = 375
$ws.Range("N71").Value = "`$/paquete 3 unidades"
$ws.Range("O71").Value = "Provincia del Elquí"
$ws.Range("P71").Value = 125
$ws.Range("Q71").Value = 3
$ws.Range("R71").Value = "Hortaliza"
